$wb = $excel.ActiveWorkbook

# --- credit_risk_parameters (sheet4): add recovery_rate / premium_frequency columns ---
$ws4 = $wb.Worksheets.Item("credit_risk_parameters")

$ws4.Range("E1").Value = "recovery_rate"
$ws4.Range("F1").Value = "premium_frequency"
$ws4.Range("E2").Value = 0.4
$ws4.Range("F2").Value = 4

# Make credit_risk_parameters the active sheet / tab (activeTab moves from 1 -> 3),
# and leave the new columns selected, matching the saved selection in the sheet.
$ws4.Activate()
$ws4.Range("E1:F2").Select()
